# Refresh the 'cryptos' price table (coinranking.com snapshot).
# All Coin/Link/Price/Volume(1h) cells are plain text in the source sheet
# (t="inlineStr"), so every literal below is written as text - including
# the Price column values that look numeric - to avoid Excel silently
# reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.824.69'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '3.119.15'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = "'248.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.06%  '

$ws.Range("D6").Value = "'623.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.74%  '

$ws.Range("E7").Value = '  +7.63%  '

$ws.Range("E8").Value = '  +2.05%  '

$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("E10").Value = '  -0.71%  '

$ws.Range("D11").Value = "'0.761"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.37%  '

$ws.Range("E12").Value = '  +3.17%  '

$ws.Range("D13").Value = "'0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.42%  '

$ws.Range("D14").Value = "'35.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.59%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '91.681.67'
$ws.Range("E15").Value = '  +0.98%  '

$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").Value = "'5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.47%  '

$ws.Range("D17").Value = '3.702.86'
$ws.Range("E17").Value = '  -0.60%  '

$ws.Range("D18").Value = '3.099.67'
$ws.Range("E18").Value = '  -2.04%  '

$ws.Range("D19").Value = "'3.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").Value = "'14.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.96%  '

$ws.Range("E21").Value = '  +1.89%  '

$ws.Range("D22").Value = "'5.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.28%  '

$ws.Range("D23").Value = "'447.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.53%  '

$ws.Range("D24").Value = "'9.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.84%  '

$ws.Range("D25").Value = "'5.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.18%  '

$ws.Range("D26").Value = "'91.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.46%  '

$ws.Range("D27").Value = "'12.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.61%  '

$ws.Range("D28").Value = '3.268.24'
$ws.Range("E28").Value = '  -0.84%  '

$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("E30").Value = '  +16.49%  '

$ws.Range("D31").Value = "'0.238"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +22.79%  '

$ws.Range("D32").Value = "'9.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.46%  '

$ws.Range("E33").Value = '  +16.83%  '

$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.67%  '

$ws.Range("E35").Value = '  +31.13%  '

$ws.Range("D36").Value = "'7.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.26%  '

$ws.Range("D37").Value = "'26.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.84%  '

$ws.Range("D38").Value = "'4.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +26.49%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = "'3.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.39%  '

$ws.Range("B40").Value = 'PancakeSwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D40").Value = "'1.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.53%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = "'498.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.98%  '

$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("E43").Value = '  +1.89%  '

$ws.Range("D44").Value = "'22.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("E46").Value = '  -0.38%  '

$ws.Range("D47").Value = "'0.700"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.10%  '

$ws.Range("D48").Value = "'153.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.91%  '

$ws.Range("D49").Value = "'4.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.67%  '

$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").Value = "'44.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.46%  '
